# v1.0.14 updates - Flow Group datasheet refresh
# - unhide Reference/Data columns C & D and widen column A to fit new names
# - append new Flow Group entries (LULC / Net Biome-Ecosystem-Primary Productivity / Q10 / Emission Rh)
# - re-sort the Name column alphabetically
# - clear the leftover per-cell formatting on the data rows (left as default/no style)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Group")

# Final (alphabetically sorted) list of Name / Description pairs for the data rows.
$rows = @(
    @("Biomass Turnover: Total", ""),
    @("Decay: Total", ""),
    @("Emission: Total Rh", "Total carbon emissions from all ecosystem components (DOM+ Biomass)"),
    @("LULC: Emission", ""),
    @("LULC: Emission CH4", ""),
    @("LULC: Emission CO", ""),
    @("LULC: Emission CO2", ""),
    @("LULC: Emission DOM", ""),
    @("LULC: Emission Live", ""),
    @("LULC: Harvest", ""),
    @("LULC: Mortality", ""),
    @("LULC: Transfer", ""),
    @("Net Biome Productivity", ""),
    @("Net Ecosystem Productivity", ""),
    @("Net Growth: Total", "Net biomass increment before losses from disturbances"),
    @("Net Primary Productivity", ""),
    @("Q10 Fast Flows", ""),
    @("Q10 Slow Flows", ""),
    @("Transfer: Total", "")
)

$startRow = 2
$endRow = $startRow + $rows.Count - 1

# Wipe the old (pre-sort) contents first so no stale leftover values survive
# in columns that end up blank for a given row after the re-sort/append.
$ws.Range("A" + $startRow + ":B" + $endRow).ClearContents()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $name = $rows[$i][0]
    $desc = $rows[$i][1]
    $ws.Cells.Item($r, 1).Value = $name
    if ($desc -ne "") {
        $ws.Cells.Item($r, 2).Value = $desc
    }
}

# Drop the bold/left-aligned leftover cell formatting on the data block so the
# cells fall back to the default (unstyled) look - matches the refreshed rows.
$dataRange = $ws.Range("A" + $startRow + ":B" + $endRow)
$dataRange.Style = "Normal"

# Reference/Data-location helper columns are no longer hidden.
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(4).Hidden = $false

# Column A needs to be wider to fit the new (longer) flow-group names.
$ws.Columns.Item(1).AutoFit()

# Leave the selection where the editor ended up.
$ws.Activate()
$ws.Range("B19").Select()
